$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.02703937269998
$ws.Range("D2").Value = 1.031234967875891
$ws.Range("E2").Value = 1.027188073291826
$ws.Range("F2").Value = 1.038153356346304
$ws.Range("I2").Value = 1.034242005509201
$ws.Range("J2").Value = 1.032199465073358
$ws.Range("K2").Value = 1.03404367325732
$ws.Range("L2").Value = 1.030008531096877
$ws.Range("M2").Value = 1.040942200810303
$ws.Range("N2").Value = 1.033665306708703

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.027848129690753
$ws.Range("D3").Value = 1.031834975261162
$ws.Range("E3").Value = 1.027871260608862
$ws.Range("F3").Value = 1.040298366781818
$ws.Range("I3").Value = 1.034463648184063
$ws.Range("J3").Value = 1.032648961335179
$ws.Range("K3").Value = 1.034452761405802
$ws.Range("L3").Value = 1.030499727670097
$ws.Range("M3").Value = 1.042893641293054
$ws.Range("N3").Value = 1.034115441306766

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.028371808347861
$ws.Range("D4").Value = 1.032223424316707
$ws.Range("E4").Value = 1.028314021422015
$ws.Range("F4").Value = 1.041681326512331
$ws.Range("I4").Value = 1.034605753672601
$ws.Range("J4").Value = 1.032939513661906
$ws.Range("K4").Value = 1.034716977254562
$ws.Range("L4").Value = 1.030817579836671
$ws.Range("M4").Value = 1.044151008114646
$ws.Range("N4").Value = 1.034406406251124

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.028592048155122
$ws.Range("D5").Value = 1.032386776657444
$ws.Range("E5").Value = 1.028500323056785
$ws.Range("F5").Value = 1.042261553888203
$ws.Range("I5").Value = 1.034665181594301
$ws.Range("J5").Value = 1.033061589389469
$ws.Range("K5").Value = 1.034827935896979
$ws.Range("L5").Value = 1.030951208054206
$ws.Range("M5").Value = 1.044678352718679
$ws.Range("N5").Value = 1.034528655340222

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.028629032372199
$ws.Range("D6").Value = 1.032414207054509
$ws.Range("E6").Value = 1.028531613569529
$ws.Range("F6").Value = 1.042358908920885
$ws.Range("I6").Value = 1.034675141462013
$ws.Range("J6").Value = 1.033082082184986
$ws.Range("K6").Value = 1.034846559438197
$ws.Range("L6").Value = 1.030973644992032
$ws.Range("M6").Value = 1.044766823544397
$ws.Range("N6").Value = 1.03454917723786

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.028374750868556
$ws.Range("D7").Value = 1.032225606849363
$ws.Range("E7").Value = 1.028316510146452
$ws.Range("F7").Value = 1.041689084096062
$ws.Range("I7").Value = 1.034606548981077
$ws.Range("J7").Value = 1.032941145128818
$ws.Range("K7").Value = 1.034718460351927
$ws.Range("L7").Value = 1.030819365371618
$ws.Range("M7").Value = 1.044158059406887
$ws.Range("N7").Value = 1.034408040034906

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.027312621263837
$ws.Range("D8").Value = 1.031437700604617
$ws.Range("E8").Value = 1.027418815630837
$ws.Range("F8").Value = 1.038879326896435
$ws.Range("I8").Value = 1.034317182927313
$ws.Range("J8").Value = 1.032351436806666
$ws.Range("K8").Value = 1.034182028279333
$ws.Range("L8").Value = 1.030174529873803
$ws.Range("M8").Value = 1.041602822679675
$ws.Range("N8").Value = 1.033817494259297

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.025443786503424
$ws.Range("D9").Value = 1.030050894204425
$ws.Range("E9").Value = 1.025842309018219
$ws.Range("F9").Value = 1.03388850818231
$ws.Range("I9").Value = 1.03379719186911
$ws.Range("J9").Value = 1.031309992310826
$ws.Range("K9").Value = 1.03323300164715
$ws.Range("L9").Value = 1.02903838247585
$ws.Range("M9").Value = 1.037057989337753
$ws.Range("N9").Value = 1.032774570792816

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.024199792276826
$ws.Range("D10").Value = 1.02912744817487
$ws.Range("E10").Value = 1.02479494876841
$ws.Range("F10").Value = 1.030532748787937
$ws.Range("I10").Value = 1.033443685973224
$ws.Range("J10").Value = 1.030614156143487
$ws.Range("K10").Value = 1.032597782623477
$ws.Range("L10").Value = 1.028281065250112
$ws.Range("M10").Value = 1.03399799561269
$ws.Range("N10").Value = 1.032077746458306

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.023661584204382
$ws.Range("D11").Value = 1.028727849231383
$ws.Range("E11").Value = 1.024342304368504
$ws.Range("F11").Value = 1.029072475146902
$ws.Range("I11").Value = 1.033288976266557
$ws.Range("J11").Value = 1.030312487119152
$ws.Range("K11").Value = 1.032322122889462
$ws.Range("L11").Value = 1.027953169867739
$ws.Range("M11").Value = 1.03266545161961
$ws.Range("N11").Value = 1.031775649029355

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.023461737665005
$ws.Range("D12").Value = 1.028579459616259
$ws.Range("E12").Value = 1.024174303446311
$ws.Range("F12").Value = 1.028528945022689
$ws.Range("I12").Value = 1.033231262609389
$ws.Range("J12").Value = 1.030200378605753
$ws.Range("K12").Value = 1.032219639349856
$ws.Range("L12").Value = 1.027831379392985
$ws.Range("M12").Value = 1.032169317597962
$ws.Range("N12").Value = 1.031663381309008

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.023504602347998
$ws.Range("D13").Value = 1.028611287913399
$ws.Range("E13").Value = 1.024210334258554
$ws.Range("F13").Value = 1.028645585332038
$ws.Range("I13").Value = 1.033243653609569
$ws.Range("J13").Value = 1.030224428761067
$ws.Range("K13").Value = 1.032241626549958
$ws.Range("L13").Value = 1.02785750365692
$ws.Range("M13").Value = 1.032275793430473
$ws.Range("N13").Value = 1.031687465618301

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.023645063439967
$ws.Range("D14").Value = 1.028715582497882
$ws.Range("E14").Value = 1.024328414675995
$ws.Range("F14").Value = 1.029027569828178
$ws.Range("I14").Value = 1.033284210695887
$ws.Range("J14").Value = 1.030303221320803
$ws.Range("K14").Value = 1.032313653430876
$ws.Range("L14").Value = 1.027943102536456
$ws.Range("M14").Value = 1.032624465049578
$ws.Range("N14").Value = 1.03176637007251

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.023731615162326
$ws.Range("D15").Value = 1.028779847038752
$ws.Range("E15").Value = 1.024401185346453
$ws.Range("F15").Value = 1.029262773636501
$ws.Range("I15").Value = 1.033309166408807
$ws.Range("J15").Value = 1.030351760701244
$ws.Range("K15").Value = 1.032358019471481
$ws.Range("L15").Value = 1.027995843431162
$ws.Range("M15").Value = 1.032839137361026
$ws.Range("N15").Value = 1.031814978384439

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.024235520877743
$ws.Range("D16").Value = 1.029153973727241
$ws.Range("E16").Value = 1.024825007668659
$ws.Range("F16").Value = 1.030629507372198
$ws.Range("I16").Value = 1.033453918898591
$ws.Range("J16").Value = 1.030634169182465
$ws.Range("K16").Value = 1.032616064470741
$ws.Range("L16").Value = 1.028302827199409
$ws.Range("M16").Value = 1.034086270277056
$ws.Range("N16").Value = 1.032097787918095

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.024551728254048
$ws.Range("D17").Value = 1.029388723131685
$ws.Range("E17").Value = 1.025091093526089
$ws.Range("F17").Value = 1.031484867985168
$ws.Range("I17").Value = 1.033544278514675
$ws.Range("J17").Value = 1.030811218203818
$ws.Range("K17").Value = 1.032777766955242
$ws.Range("L17").Value = 1.028495397691258
$ws.Range("M17").Value = 1.034866519513109
$ws.Range("N17").Value = 1.032275088369367

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.024736210029209
$ws.Range("D18").Value = 1.029525673488932
$ws.Range("E18").Value = 1.02524638075227
$ws.Range("F18").Value = 1.031983092424815
$ws.Range("I18").Value = 1.033596825668775
$ws.Range("J18").Value = 1.030914452459356
$ws.Range("K18").Value = 1.032872026855303
$ws.Range("L18").Value = 1.028607723484214
$ws.Range("M18").Value = 1.035320899891545
$ws.Range("N18").Value = 1.032378469229389

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.024799120890694
$ws.Range("D19").Value = 1.029572374253252
$ws.Range("E19").Value = 1.025299343882507
$ws.Range("F19").Value = 1.032152857648557
$ws.Range("I19").Value = 1.033614716115589
$ws.Range("J19").Value = 1.030949646660626
$ws.Range("K19").Value = 1.032904157142307
$ws.Range("L19").Value = 1.028646024123215
$ws.Range("M19").Value = 1.0354757097728
$ws.Range("N19").Value = 1.032413713410461

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.024517797707558
$ws.Range("D20").Value = 1.029363534154337
$ws.Range("E20").Value = 1.025062536353888
$ws.Range("F20").Value = 1.031393167823111
$ws.Range("I20").Value = 1.033534600143443
$ws.Range("J20").Value = 1.0307922261856
$ws.Range("K20").Value = 1.032760423853766
$ws.Range("L20").Value = 1.028474736420557
$ws.Range("M20").Value = 1.034782881472216
$ws.Range("N20").Value = 1.032256069380305

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.023603699272149
$ws.Range("D21").Value = 1.028684869257411
$ws.Range("E21").Value = 1.024293639304166
$ws.Range("F21").Value = 1.028915116114729
$ws.Range("I21").Value = 1.033272274489901
$ws.Range("J21").Value = 1.030280020399982
$ws.Range("K21").Value = 1.032292445830417
$ws.Range("L21").Value = 1.027917895674503
$ws.Range("M21").Value = 1.032521822463
$ws.Range("N21").Value = 1.03174313620372

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.023029362323834
$ws.Range("D22").Value = 1.028258392753692
$ws.Range("E22").Value = 1.023810963574173
$ws.Range("F22").Value = 1.027350571139495
$ws.Range("I22").Value = 1.033105906832217
$ws.Range("J22").Value = 1.02995765686524
$ws.Range("K22").Value = 1.031997681444312
$ws.Range("L22").Value = 1.02756781377343
$ws.Range("M22").Value = 1.031093430805237
$ws.Range("N22").Value = 1.031420314875783

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.02333379184215
$ws.Range("D23").Value = 1.028484454412749
$ws.Range("E23").Value = 1.024066766800054
$ws.Range("F23").Value = 1.028180594111939
$ws.Range("I23").Value = 1.033194237757805
$ws.Range("J23").Value = 1.030128578162321
$ws.Range("K23").Value = 1.03215399176043
$ws.Range("L23").Value = 1.027753396266891
$ws.Range("M23").Value = 1.031851301881057
$ws.Range("N23").Value = 1.031591478900711

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.024533129327482
$ws.Range("D24").Value = 1.029374915892058
$ws.Range("E24").Value = 1.025075439852111
$ws.Range("F24").Value = 1.031434605318949
$ws.Range("I24").Value = 1.03353897387135
$ws.Range("J24").Value = 1.030800807970918
$ws.Range("K24").Value = 1.032768260635561
$ws.Range("L24").Value = 1.028484072351444
$ws.Range("M24").Value = 1.03482067614418
$ws.Range("N24").Value = 1.032264663352742

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.0259265931655
$ws.Range("D25").Value = 1.030409226298304
$ws.Range("E25").Value = 1.026249235035168
$ws.Range("F25").Value = 1.035183646679164
$ws.Range("I25").Value = 1.033932823998092
$ws.Range("J25").Value = 1.031579502639784
$ws.Range("K25").Value = 1.033478794187562
$ws.Range("L25").Value = 1.029332085303442
$ws.Range("M25").Value = 1.038238114246747
$ws.Range("N25").Value = 1.033044463857355
